# NatmiData TPM update: Nppc-Npr3.xlsx (YoungD4)
# The underlying NATMI computation was rerun with updated TPM values. The sending/target
# cluster matrix now spans all three clusters (ECs, FAPs, MuSCs) rather than just two,
# so rows 2-7 are refreshed in place and rows 8-10 are brand-new additions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs (sending) -> ECs (target), Nppc/Npr3
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Nppc"
$ws.Range("C2").Value = "Npr3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.08453933333333334
$ws.Range("H2").Value = 0.253618
$ws.Range("I2").Value = 0.04191805799619459
$ws.Range("J2").Value = 0.04191805799619459
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04090066666666667
$ws.Range("N2").Value = 0.122702
$ws.Range("O2").Value = 0.09164231251535751
$ws.Range("P2").Value = 0.0916423125153575
$ws.Range("Q2").Value = 0.003457715092888889
$ws.Range("R2").Value = 0.031119435836
$ws.Range("S2").Value = 0.003841467770924146
$ws.Range("T2").Value = 0.003841467770924145

# Row 3: ECs (sending) -> FAPs (target), Nppc/Npr3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Nppc"
$ws.Range("C3").Value = "Npr3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.08453933333333334
$ws.Range("H3").Value = 0.253618
$ws.Range("I3").Value = 0.04191805799619459
$ws.Range("J3").Value = 0.04191805799619459
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3809676666666666
$ws.Range("N3").Value = 1.142903
$ws.Range("O3").Value = 0.8535987506376395
$ws.Range("P3").Value = 0.8535987506376393
$ws.Range("Q3").Value = 0.03220675256155556
$ws.Range("R3").Value = 0.289860773054
$ws.Range("S3").Value = 0.03578120193470782
$ws.Range("T3").Value = 0.03578120193470781

# Row 4: ECs (sending) -> MuSCs (target), Nppc/Npr3
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Nppc"
$ws.Range("C4").Value = "Npr3"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.08453933333333334
$ws.Range("H4").Value = 0.253618
$ws.Range("I4").Value = 0.04191805799619459
$ws.Range("J4").Value = 0.04191805799619459
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02443933333333333
$ws.Range("N4").Value = 0.073318
$ws.Range("O4").Value = 0.05475893684700315
$ws.Range("P4").Value = 0.05475893684700314
$ws.Range("Q4").Value = 0.002066084947111111
$ws.Range("R4").Value = 0.018594764524
$ws.Range("S4").Value = 0.002295388290562635
$ws.Range("T4").Value = 0.002295388290562635

# Row 5: FAPs (sending) -> ECs (target), Nppc/Npr3
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Nppc"
$ws.Range("C5").Value = "Npr3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.435073666666667
$ws.Range("H5").Value = 4.305221
$ws.Range("I5").Value = 0.7115681992777912
$ws.Range("J5").Value = 0.7115681992777912
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.04090066666666667
$ws.Range("N5").Value = 0.122702
$ws.Range("O5").Value = 0.09164231251535751
$ws.Range("P5").Value = 0.0916423125153575
$ws.Range("Q5").Value = 0.05869546968244444
$ws.Range("R5").Value = 0.5282592271419999
$ws.Range("S5").Value = 0.06520975529420553
$ws.Range("T5").Value = 0.06520975529420553

# Row 6: FAPs (sending) -> FAPs (target), Nppc/Npr3
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Nppc"
$ws.Range("C6").Value = "Npr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.435073666666667
$ws.Range("H6").Value = 4.305221
$ws.Range("I6").Value = 0.7115681992777912
$ws.Range("J6").Value = 0.7115681992777912
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.3809676666666666
$ws.Range("N6").Value = 1.142903
$ws.Range("O6").Value = 0.8535987506376395
$ws.Range("P6").Value = 0.8535987506376393
$ws.Range("Q6").Value = 0.5467166662847778
$ws.Range("R6").Value = 4.920449996563
$ws.Range("S6").Value = 0.6073937258969975
$ws.Range("T6").Value = 0.6073937258969974

# Row 7: FAPs (sending) -> MuSCs (target), Nppc/Npr3
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Nppc"
$ws.Range("C7").Value = "Npr3"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.435073666666667
$ws.Range("H7").Value = 4.305221
$ws.Range("I7").Value = 0.7115681992777912
$ws.Range("J7").Value = 0.7115681992777912
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.02443933333333333
$ws.Range("N7").Value = 0.073318
$ws.Range("O7").Value = 0.05475893684700315
$ws.Range("P7").Value = 0.05475893684700314
$ws.Range("Q7").Value = 0.03507224369755555
$ws.Range("R7").Value = 0.3156501932779999
$ws.Range("S7").Value = 0.03896471808658832
$ws.Range("T7").Value = 0.03896471808658832

# Row 8: MuSCs (sending) -> ECs (target), Nppc/Npr3
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Nppc"
$ws.Range("C8").Value = "Npr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.497163
$ws.Range("H8").Value = 1.491489
$ws.Range("I8").Value = 0.2465137427260142
$ws.Range("J8").Value = 0.2465137427260142
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.04090066666666667
$ws.Range("N8").Value = 0.122702
$ws.Range("O8").Value = 0.09164231251535751
$ws.Range("P8").Value = 0.0916423125153575
$ws.Range("Q8").Value = 0.020334298142
$ws.Range("R8").Value = 0.183008683278
$ws.Range("S8").Value = 0.02259108945022784
$ws.Range("T8").Value = 0.02259108945022783

# Row 9: MuSCs (sending) -> FAPs (target), Nppc/Npr3
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Nppc"
$ws.Range("C9").Value = "Npr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.497163
$ws.Range("H9").Value = 1.491489
$ws.Range("I9").Value = 0.2465137427260142
$ws.Range("J9").Value = 0.2465137427260142
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.3809676666666666
$ws.Range("N9").Value = 1.142903
$ws.Range("O9").Value = 0.8535987506376395
$ws.Range("P9").Value = 0.8535987506376393
$ws.Range("Q9").Value = 0.189403028063
$ws.Range("R9").Value = 1.704627252567
$ws.Range("S9").Value = 0.2104238228059342
$ws.Range("T9").Value = 0.2104238228059342

# Row 10: MuSCs (sending) -> MuSCs (target), Nppc/Npr3
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Nppc"
$ws.Range("C10").Value = "Npr3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.497163
$ws.Range("H10").Value = 1.491489
$ws.Range("I10").Value = 0.2465137427260142
$ws.Range("J10").Value = 0.2465137427260142
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.02443933333333333
$ws.Range("N10").Value = 0.073318
$ws.Range("O10").Value = 0.05475893684700315
$ws.Range("P10").Value = 0.05475893684700314
$ws.Range("Q10").Value = 0.012150332278
$ws.Range("R10").Value = 0.109352990502
$ws.Range("S10").Value = 0.0134988304698522
$ws.Range("T10").Value = 0.01349883046985219
